$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 40 and 41 swap places: Maker moves up to row 40, FirstDigitalUSD moves
# down to row 41. Coin name / link / price / volume are all updated.
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "'3.073.82"
$ws.Range("E40").Value = "'  -5.20%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = "'  -0.14%  "

# Remaining Price / Volume(1h) cell updates (GitHub Actions refresh of the
# cryptos list). Values are forced to text with a leading apostrophe so
# Excel doesn't reinterpret thousands-separator dots / trailing zeros as
# numbers.
$ws.Range("D2").Value = "'67.309.88"
$ws.Range("E2").Value = "'  -1.77%  "
$ws.Range("D3").Value = "'3.507.71"
$ws.Range("E3").Value = "'  -3.78%  "
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("D5").Value = "'201.92"
$ws.Range("E5").Value = "'  +2.68%  "
$ws.Range("D6").Value = "'553.31"
$ws.Range("E6").Value = "'  -5.14%  "
$ws.Range("D7").Value = "'3.495.14"
$ws.Range("E7").Value = "'  -3.99%  "
$ws.Range("D8").Value = "'0.608"
$ws.Range("E8").Value = "'  -2.08%  "
$ws.Range("E9").Value = "'  -0.17%  "
$ws.Range("D10").Value = "'0.655"
$ws.Range("E10").Value = "'  -3.90%  "
$ws.Range("D11").Value = "'62.92"
$ws.Range("E11").Value = "'  +10.93%  "
$ws.Range("E12").Value = "'  -7.17%  "
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("E13").Value = "'  -7.74%  "
$ws.Range("D14").Value = "'9.84"
$ws.Range("E14").Value = "'  -3.34%  "
$ws.Range("D15").Value = "'4.060.70"
$ws.Range("E15").Value = "'  -4.07%  "
$ws.Range("D16").Value = "'3.499.55"
$ws.Range("E16").Value = "'  -4.08%  "
$ws.Range("E17").Value = "'  -1.92%  "
$ws.Range("D18").Value = "'18.43"
$ws.Range("E18").Value = "'  -1.32%  "
$ws.Range("D19").Value = "'66.992.47"
$ws.Range("E19").Value = "'  -2.19%  "
$ws.Range("D20").Value = "'11.83"
$ws.Range("E20").Value = "'  -6.06%  "
$ws.Range("D21").Value = "'1.03"
$ws.Range("E21").Value = "'  -5.77%  "
$ws.Range("D22").Value = "'392.07"
$ws.Range("E22").Value = "'  -2.73%  "
$ws.Range("D23").Value = "'12.56"
$ws.Range("E23").Value = "'  -3.48%  "
$ws.Range("D24").Value = "'4.00"
$ws.Range("E24").Value = "'  -6.43%  "
$ws.Range("D25").Value = "'82.88"
$ws.Range("E25").Value = "'  -3.73%  "
$ws.Range("D26").Value = "'3.92"
$ws.Range("E26").Value = "'  +1.36%  "
$ws.Range("D27").Value = "'12.30"
$ws.Range("E27").Value = "'  -2.87%  "
$ws.Range("D28").Value = "'2.82"
$ws.Range("E28").Value = "'  -4.96%  "
$ws.Range("D29").Value = "'8.84"
$ws.Range("E29").Value = "'  -4.00%  "
$ws.Range("D30").Value = "'31.04"
$ws.Range("E30").Value = "'  -2.50%  "
$ws.Range("D31").Value = "'678.98"
$ws.Range("E31").Value = "'  -3.65%  "
$ws.Range("D32").Value = "'7.00"
$ws.Range("E32").Value = "'  -13.97%  "
$ws.Range("D33").Value = "'11.75"
$ws.Range("E33").Value = "'  -4.16%  "
$ws.Range("D34").Value = "'63.71"
$ws.Range("E34").Value = "'  -1.96%  "
$ws.Range("E35").Value = "'  -6.39%  "
$ws.Range("D36").Value = "'38.74"
$ws.Range("E36").Value = "'  -9.37%  "
$ws.Range("E37").Value = "'  +0.06%  "
$ws.Range("D38").Value = "'0.397"
$ws.Range("E38").Value = "'  -7.57%  "
$ws.Range("D39").Value = "'0.131"
$ws.Range("E39").Value = "'  -5.07%  "
$ws.Range("D42").Value = "'2.99"
$ws.Range("E42").Value = "'  -4.29%  "
$ws.Range("D43").Value = "'2.60"
$ws.Range("E43").Value = "'  -9.52%  "
$ws.Range("D44").Value = "'0.0₃0676"
$ws.Range("E44").Value = "'  -14.48%  "
$ws.Range("E45").Value = "'  +5.49%  "
$ws.Range("D46").Value = "'2.72"
$ws.Range("E46").Value = "'  -11.48%  "
$ws.Range("D47").Value = "'0.0401"
$ws.Range("E47").Value = "'  -5.27%  "
$ws.Range("E48").Value = "'  -3.85%  "
$ws.Range("D49").Value = "'137.73"
$ws.Range("E49").Value = "'  -3.63%  "
$ws.Range("D50").Value = "'8.21"
$ws.Range("E50").Value = "'  -7.98%  "
$ws.Range("D51").Value = "'2.85"
$ws.Range("E51").Value = "'  -8.59%  "
